$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")
$summary = $wb.Worksheets.Item("Summary")

# Columns: D=MyForecast, G=Trend, H=Inventory Coverage, I=Stockout Risk,
#          J=Reorder Urgency, L=Seasonality Index
# "added 4wk low sales check": trend flips to "Low Volume Season" for the
# whole series, forecast/coverage/seasonality numbers are recomputed, and
# several weeks at the tail drop to a zero forecast with relaxed urgency.

$ws.Range("D2").Value = 112
$ws.Range("G2").Value = "Low Volume Season"
$ws.Range("H2").Value = 8.46
$ws.Range("L2").Value = 1.05

$ws.Range("D3").Value = 120
$ws.Range("G3").Value = "Low Volume Season"
$ws.Range("H3").Value = 6.97
$ws.Range("L3").Value = 1.13

$ws.Range("D4").Value = 81
$ws.Range("G4").Value = "Low Volume Season"
$ws.Range("H4").Value = 8.84
$ws.Range("L4").Value = 0.92

$ws.Range("D5").Value = 29
$ws.Range("G5").Value = "Low Volume Season"
$ws.Range("H5").Value = 21.9
$ws.Range("L5").Value = 0.98

$ws.Range("D6").Value = 11
$ws.Range("G6").Value = "Low Volume Season"
$ws.Range("H6").Value = 55.09
$ws.Range("L6").Value = 0.9

$ws.Range("D7").Value = 27
$ws.Range("G7").Value = "Low Volume Season"
$ws.Range("H7").Value = 22.04
$ws.Range("L7").Value = 0.9

$ws.Range("D8").Value = 13
$ws.Range("G8").Value = "Low Volume Season"
$ws.Range("H8").Value = 43.69
$ws.Range("L8").Value = 1.19

$ws.Range("D9").Value = 0
$ws.Range("G9").Value = "Low Volume Season"
$ws.Range("H9").Value = ""
$ws.Range("L9").Value = 1.09

$ws.Range("D10").Value = 0
$ws.Range("G10").Value = "Low Volume Season"
$ws.Range("H10").Value = ""
$ws.Range("L10").Value = 1.2

$ws.Range("D11").Value = 0
$ws.Range("G11").Value = "Low Volume Season"
$ws.Range("H11").Value = ""
$ws.Range("L11").Value = 0.97

$ws.Range("D12").Value = 37
$ws.Range("G12").Value = "Low Volume Season"
$ws.Range("H12").Value = 15
$ws.Range("L12").Value = 1.17

$ws.Range("D13").Value = 72
$ws.Range("G13").Value = "Low Volume Season"
$ws.Range("H13").Value = 7.19
$ws.Range("L13").Value = 0.84

$ws.Range("D14").Value = 0
$ws.Range("G14").Value = "Low Volume Season"
$ws.Range("H14").Value = ""
$ws.Range("J14").Value = "Normal"
$ws.Range("L14").Value = 1.01

$ws.Range("D15").Value = 0
$ws.Range("G15").Value = "Low Volume Season"
$ws.Range("H15").Value = ""
$ws.Range("I15").Value = "Low"
$ws.Range("J15").Value = "Normal"
$ws.Range("L15").Value = 0.91

$ws.Range("D16").Value = 0
$ws.Range("G16").Value = "Low Volume Season"
$ws.Range("H16").Value = ""
$ws.Range("I16").Value = "Low"
$ws.Range("J16").Value = "Normal"
$ws.Range("L16").Value = 0.82

$ws.Range("D17").Value = 0
$ws.Range("G17").Value = "Low Volume Season"
$ws.Range("H17").Value = ""
$ws.Range("I17").Value = "Low"
$ws.Range("J17").Value = "Normal"
$ws.Range("L17").Value = 0.93

# Summary sheet totals recomputed off the new MyForecast column.
# These cells are stored as text in the workbook, so force text with a
# leading quote (same as typing '502 into the cell in Excel) rather than
# letting AutoCorrect turn the digits into a number.
$summary.Range("B9").Value = "'502"
$summary.Range("B10").Value = "'393"
$summary.Range("B11").Value = "'342"
$summary.Range("B12").Value = "'120"
$summary.Range("B14").Value = "'0"
